$wb = $excel.ActiveWorkbook

# --- Sheet1: update the Overview description text and move selection to F2 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Value = "This API is used to retrieve number of leads for an Agent. 11"
$ws1.Activate()
$ws1.Range("F2").Select()

# --- Sheet3: update the Response sample JSON and move selection to A2 ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A2").Value = '{ "leadCount": { "leadCount": "11" }, "systemInformation": { "errorCode": "", "errorMessage": "", "timestamp": "", "details": "333" } }'
$ws3.Activate()
$ws3.Range("A2").Select()
